# Auto-generated edit script
# Updates cached market-data columns (H-N) on each leve-profit sheet
# to match the refreshed values from the scheduled Kujata data-pull.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 6600.5
$ws.Range("I28").Value = 7007.2
$ws.Range("K28").Value = 7007.2
$ws.Range("M28").Value = -6522.2
$ws.Range("H76").Value = 4999.8887
$ws.Range("I76").Value = 5333
$ws.Range("J76").Value = 4833.3335
$ws.Range("K76").Value = 5333
$ws.Range("L76").Value = 4833.3335
$ws.Range("M76").Value = -5018
$ws.Range("N76").Value = -5463.3335
$ws.Range("H79").Value = 4999.8887
$ws.Range("I79").Value = 5333
$ws.Range("J79").Value = 4833.3335
$ws.Range("K79").Value = 5333
$ws.Range("L79").Value = 4833.3335
$ws.Range("M79").Value = -4241
$ws.Range("N79").Value = -7017.3335
$ws.Range("H111").Value = 1885.1111
$ws.Range("I111").Value = 1357.5454
$ws.Range("J111").Value = 2714.1428
$ws.Range("K111").Value = 4072.6362
$ws.Range("L111").Value = 8142.428400000001
$ws.Range("M111").Value = -1005.6362
$ws.Range("N111").Value = -14276.4284
$ws.Range("H137").Value = 2000.4
$ws.Range("I137").Value = 1333.6666
$ws.Range("J137").Value = 2444.889
$ws.Range("K137").Value = 4000.9998
$ws.Range("L137").Value = 7334.667
$ws.Range("M137").Value = -1450.9998
$ws.Range("N137").Value = -12434.667

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2386.111
$ws.Range("I2").Value = 1072.6666
$ws.Range("J2").Value = 5013
$ws.Range("K2").Value = 1072.6666
$ws.Range("L2").Value = 5013
$ws.Range("M2").Value = -959.6666
$ws.Range("N2").Value = -5239
$ws.Range("H62").Value = 75000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 75000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 75000
$ws.Range("M62").Value = $null
$ws.Range("N62").Value = -76248
$ws.Range("H65").Value = 75000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 75000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 225000
$ws.Range("M65").Value = $null
$ws.Range("N65").Value = -231240
$ws.Range("H94").Value = 49999.5
$ws.Range("J94").Value = 49999.5
$ws.Range("L94").Value = 49999.5
$ws.Range("N94").Value = -51801.5
$ws.Range("H110").Value = 1728.875
$ws.Range("I110").Value = 1384.5385
$ws.Range("K110").Value = 1384.5385
$ws.Range("M110").Value = 660.4614999999999
$ws.Range("H116").Value = 2386.111
$ws.Range("I116").Value = 1072.6666
$ws.Range("J116").Value = 5013
$ws.Range("K116").Value = 1072.6666
$ws.Range("L116").Value = 5013
$ws.Range("M116").Value = 1221.3334
$ws.Range("N116").Value = -9601

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2386.111
$ws.Range("I3").Value = 1072.6666
$ws.Range("J3").Value = 5013
$ws.Range("K3").Value = 1072.6666
$ws.Range("L3").Value = 5013
$ws.Range("M3").Value = -958.6666
$ws.Range("N3").Value = -5241
$ws.Range("H86").Value = 2939.5151
$ws.Range("I86").Value = 2946.4814
$ws.Range("K86").Value = 2946.4814
$ws.Range("M86").Value = -1823.4814
$ws.Range("H89").Value = 2939.5151
$ws.Range("I89").Value = 2946.4814
$ws.Range("K89").Value = 14732.407
$ws.Range("M89").Value = -9116.407000000001
$ws.Range("H94").Value = 15625759
$ws.Range("I94").Value = 25000590
$ws.Range("J94").Value = 1039.8334
$ws.Range("K94").Value = 25000590
$ws.Range("L94").Value = 1039.8334
$ws.Range("M94").Value = -25000139
$ws.Range("N94").Value = -1941.8334
$ws.Range("H107").Value = 1273.3914
$ws.Range("I107").Value = 1035.5264
$ws.Range("K107").Value = 1035.5264
$ws.Range("M107").Value = 884.4736

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1199.4182
$ws.Range("I31").Value = 869.7273
$ws.Range("J31").Value = 1419.2122
$ws.Range("K31").Value = 869.7273
$ws.Range("L31").Value = 1419.2122
$ws.Range("M31").Value = -574.7273
$ws.Range("N31").Value = -2009.2122
$ws.Range("H34").Value = 1199.4182
$ws.Range("I34").Value = 869.7273
$ws.Range("J34").Value = 1419.2122
$ws.Range("K34").Value = 869.7273
$ws.Range("L34").Value = 1419.2122
$ws.Range("M34").Value = -667.7273
$ws.Range("N34").Value = -1823.2122
$ws.Range("H134").Value = 1467.8387
$ws.Range("I134").Value = 1457.9048
$ws.Range("J134").Value = 1488.7
$ws.Range("K134").Value = 4373.7144
$ws.Range("L134").Value = 4466.1
$ws.Range("M134").Value = -1838.7144
$ws.Range("N134").Value = -9536.1

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 4933.304
$ws.Range("J107").Value = 25524.75
$ws.Range("L107").Value = 76574.25
$ws.Range("N107").Value = -80414.25
$ws.Range("H118").Value = 1000
$ws.Range("I118").Value = 1000
$ws.Range("K118").Value = 3000
$ws.Range("M118").Value = -1757
$ws.Range("H131").Value = 21278044
$ws.Range("I131").Value = 142857970
$ws.Range("J131").Value = 1557.45
$ws.Range("K131").Value = 428573910
$ws.Range("L131").Value = 4672.35
$ws.Range("M131").Value = -428568870
$ws.Range("N131").Value = -14752.35
$ws.Range("H132").Value = 1456.2727
$ws.Range("I132").Value = 799
$ws.Range("J132").Value = 1522
$ws.Range("K132").Value = 7191
$ws.Range("L132").Value = 13698
$ws.Range("M132").Value = -4661
$ws.Range("N132").Value = -18758

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 19017
$ws.Range("J95").Value = 19017
$ws.Range("L95").Value = 19017
$ws.Range("N95").Value = -24509
$ws.Range("H113").Value = 2159.7778
$ws.Range("I113").Value = 1134
$ws.Range("J113").Value = 5750
$ws.Range("K113").Value = 1134
$ws.Range("L113").Value = 5750
$ws.Range("M113").Value = 1036
$ws.Range("N113").Value = -10090

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3015.8572
$ws.Range("I61").Value = 2400.4285
$ws.Range("J61").Value = 3631.2856
$ws.Range("K61").Value = 2400.4285
$ws.Range("L61").Value = 3631.2856
$ws.Range("M61").Value = -2198.4285
$ws.Range("N61").Value = -4035.2856
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null
$ws.Range("H113").Value = 3015.8572
$ws.Range("I113").Value = 2400.4285
$ws.Range("J113").Value = 3631.2856
$ws.Range("K113").Value = 2400.4285
$ws.Range("L113").Value = 3631.2856
$ws.Range("M113").Value = -230.4285
$ws.Range("N113").Value = -7971.2856

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 8200
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 9933.333000000001
$ws.Range("K63").Value = 3000
$ws.Range("L63").Value = 9933.333000000001
$ws.Range("M63").Value = -2376
$ws.Range("N63").Value = -11181.333
$ws.Range("H66").Value = 8200
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 9933.333000000001
$ws.Range("K66").Value = 9000
$ws.Range("L66").Value = 29799.999
$ws.Range("M66").Value = -5880
$ws.Range("N66").Value = -36039.999
$ws.Range("H92").Value = 16749.834
$ws.Range("J92").Value = 16749.834
$ws.Range("L92").Value = 16749.834
$ws.Range("N92").Value = -21741.834
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492
$ws.Range("H113").Value = 529.2353000000001
$ws.Range("I113").Value = 302
$ws.Range("K113").Value = 906
$ws.Range("M113").Value = 1264
$ws.Range("H122").Value = 11306006
$ws.Range("I122").Value = 12382612
$ws.Range("K122").Value = 37147836
$ws.Range("M122").Value = -37145386
$ws.Range("H126").Value = 55556924
$ws.Range("I126").Value = 76923850
$ws.Range("J126").Value = 2919
$ws.Range("K126").Value = 230771550
$ws.Range("L126").Value = 8757
$ws.Range("M126").Value = -230769080
$ws.Range("N126").Value = -13697
$ws.Range("H139").Value = 37126
$ws.Range("J139").Value = 38543.332
$ws.Range("L139").Value = 38543.332
$ws.Range("N139").Value = -48823.332
$ws.Range("H141").Value = 43585.832
$ws.Range("J141").Value = 43585.832
$ws.Range("L141").Value = 43585.832
$ws.Range("N141").Value = -53945.832
